# minor edits lesson 4, 6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lesson 4: row 30 (iris dataset) was missing its Lesson number -> set to 4
$ws.Range("A30").Value = 4

# Lesson 4: add a new row 50 describing the gapminder dataset
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "gapminder"
$ws.Range("C50").Value = "gapminder"
$ws.Range("D50").Value = "Life expectancy, population, and GDP by country"
$ws.Range("F50").Value = "data"

# Move the view / selection to reflect where the editor ended up working
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D51").Select()
